$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9147897958755493
$ws.Range("B1").Value = 1.905754923820496
$ws.Range("C1").Value = 3.423617362976074
$ws.Range("D1").Value = 3.744511365890503
$ws.Range("E1").Value = 0.9344961047172546
